# The document has two distinct headers (default + first-page) and two
# distinct footers (default + first-page), each carrying one inline
# picture:
#   - Headers: BTec_Logo-Orange.jpg, originally named "image1.jpg"
#       -> renamed to "image2.jpg"
#   - Footers: PearsonLogo.png, originally named "image2.png"
#       -> renamed to "image1.png"
#
# InlineShape has no settable Name property (same as real Word), so the
# rename is done by converting the inline picture to a floating Shape
# (which does expose Name), setting the new name there, and converting
# it straight back to an inline picture so the layout/structure of the
# document is left untouched.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-HeaderFooterImage($range, $newName) {
    $inlineShape = $range.InlineShapes.Item(1)
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

# Headers (BTec_Logo-Orange.jpg): image1.jpg -> image2.jpg
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    Rename-HeaderFooterImage $hdr.Range "image2.jpg"
}

# Footers (PearsonLogo.png): image2.png -> image1.png
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    Rename-HeaderFooterImage $ftr.Range "image1.png"
}
